$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Week 3")
$ws4 = $wb.Worksheets.Item("Week 4")

# --- New timesheet entries for Week 4 (rows 2-5) ---

# Row 2: 1/28 (41666) 10:00 PM - 12:00 AM, 2 hrs
$ws4.Cells.Item(2,1).Value = 41666
$ws4.Cells.Item(2,1).NumberFormat = "m/d/yy"
$ws4.Cells.Item(2,2).Value = 0.91666666666666663
$ws4.Cells.Item(2,2).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(2,3).Value = 0
$ws4.Cells.Item(2,3).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(2,4).Value = "Researched how to deploy Laravel project to shared server"
$ws4.Cells.Item(2,5).Value = 2

# Row 3: 1/29 (41667) 10:00 PM - 12:30 AM, 2.5 hrs
$ws4.Cells.Item(3,1).Value = 41667
$ws4.Cells.Item(3,1).NumberFormat = "m/d/yy"
$ws4.Cells.Item(3,2).Value = 0.91666666666666663
$ws4.Cells.Item(3,2).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(3,3).Value = 0.020833333333333332
$ws4.Cells.Item(3,3).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(3,4).Value = "Researched Deployer to upload Laravel to shared server"
$ws4.Cells.Item(3,5).Value = 2.5

# Row 4: 1/30 (41668) 10:00 AM - 12:00 PM, 2 hrs
$ws4.Cells.Item(4,1).Value = 41668
$ws4.Cells.Item(4,1).NumberFormat = "m/d/yy"
$ws4.Cells.Item(4,2).Value = 0.41666666666666669
$ws4.Cells.Item(4,2).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(4,3).Value = 0.5
$ws4.Cells.Item(4,3).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(4,4).Value = "Finished project description"
$ws4.Cells.Item(4,5).Value = 2

# Row 5: 1/30 (41668) 8:00 PM - 9:00 PM, 1 hr
$ws4.Cells.Item(5,1).Value = 41668
$ws4.Cells.Item(5,1).NumberFormat = "m/d/yy"
$ws4.Cells.Item(5,2).Value = 0.83333333333333337
$ws4.Cells.Item(5,2).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(5,3).Value = 0.875
$ws4.Cells.Item(5,3).NumberFormat = "h:mm AM/PM"
$ws4.Cells.Item(5,4).Value = "Contacted SiteGround support and adjusted settings to get Laravel project deployed live"
$ws4.Cells.Item(5,4).WrapText = $true
$ws4.Rows.Item(5).RowHeight = 26
$ws4.Cells.Item(5,5).Value = 1

# Match the row-height / sheet-format touch-ups Week 4 picks up alongside
# the new entries (rows 1/20/21 lose their custom 18pt height).
$ws4.Rows.Item(1).RowHeight = 14
$ws4.Rows.Item(20).RowHeight = 13
$ws4.Rows.Item(21).RowHeight = 13

# --- Move the active tab / selection from Week 3 to Week 4 ---
$ws3.Range("A1:XFD1048576").Select()
$ws4.Activate()
$ws4.Range("E6").Select()
